# Updated full name functions
#
# Replace `<obj>.name.full(middle=’full’)` with `<obj>.name_full()`
# for property_agent, user, and person (the "Dear ..." greeting included).

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
}

$rsquo = [char]0x2019

Replace-Text "property_agent.name.full(middle=$($rsquo)full$($rsquo)) }}" "property_agent.name_full() }}"
Replace-Text "user.name.full(middle=$($rsquo)full$($rsquo)) }}" "user.name_full() }}"
Replace-Text "person.name.full(middle=$($rsquo)full$($rsquo)) }}" "person.name_full() }}"
